$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")
$ws.Rows("1:2").Delete()
$ws.Rows("1:1").Select()
